$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Cells.Item(100, 8).Value = 2401.25
$ws.Cells.Item(100, 9).Value = 2262
$ws.Cells.Item(100, 11).Value = 2262
$ws.Cells.Item(100, 13).Value = -1721
$ws.Cells.Item(129, 8).Value = 1314.2858
$ws.Cells.Item(129, 10).Value = 2500
$ws.Cells.Item(129, 12).Value = 7500
$ws.Cells.Item(129, 14).Value = -17500
$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Cells.Item(55, 8).Value = 24100
$ws.Cells.Item(55, 10).Value = 24100
$ws.Cells.Item(55, 12).Value = 24100
$ws.Cells.Item(55, 14).Value = -24730
$ws.Cells.Item(80, 8).Value = 42000
$ws.Cells.Item(80, 10).Value = 42000
$ws.Cells.Item(80, 12).Value = 42000
$ws.Cells.Item(80, 14).Value = -43996
$ws.Cells.Item(83, 8).Value = 42000
$ws.Cells.Item(83, 10).Value = 42000
$ws.Cells.Item(83, 12).Value = 126000
$ws.Cells.Item(83, 14).Value = -135984
$ws.Cells.Item(130, 8).Value = 22869.2
$ws.Cells.Item(130, 10).Value = 22869.2
$ws.Cells.Item(130, 12).Value = 22869.2
$ws.Cells.Item(130, 14).Value = -32909.2
$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Cells.Item(35, 8).Value = 40000
$ws.Cells.Item(35, 10).Value = 40000
$ws.Cells.Item(35, 12).Value = 40000
$ws.Cells.Item(35, 14).Value = -40620
$ws.Cells.Item(39, 8).Value = 16577.666
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 16577.666
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 16577.666
$ws.Cells.Item(39, 13).Value = ""
$ws.Cells.Item(39, 14).Value = -17355.666
$ws.Cells.Item(46, 8).Value = 32987.332
$ws.Cells.Item(46, 9).Value = 8059
$ws.Cells.Item(46, 10).Value = 37973
$ws.Cells.Item(46, 11).Value = 8059
$ws.Cells.Item(46, 12).Value = 37973
$ws.Cells.Item(46, 13).Value = -7761
$ws.Cells.Item(46, 14).Value = -38569
$ws.Cells.Item(82, 8).Value = 6314.25
$ws.Cells.Item(82, 9).Value = 6314.25
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 6314.25
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -5931.25
$ws.Cells.Item(82, 14).Value = ""
$ws.Cells.Item(85, 8).Value = 6314.25
$ws.Cells.Item(85, 9).Value = 6314.25
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 6314.25
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = -4988.25
$ws.Cells.Item(85, 14).Value = ""
$ws.Cells.Item(94, 8).Value = 1643.2222
$ws.Cells.Item(94, 9).Value = 2649.5
$ws.Cells.Item(94, 10).Value = 1355.7142
$ws.Cells.Item(94, 11).Value = 2649.5
$ws.Cells.Item(94, 12).Value = 1355.7142
$ws.Cells.Item(94, 13).Value = -2198.5
$ws.Cells.Item(94, 14).Value = -2257.7142
$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Cells.Item(31, 8).Value = 56495.79
$ws.Cells.Item(31, 9).Value = 2686.5715
$ws.Cells.Item(31, 11).Value = 2686.5715
$ws.Cells.Item(31, 13).Value = -2391.5715
$ws.Cells.Item(34, 8).Value = 56495.79
$ws.Cells.Item(34, 9).Value = 2686.5715
$ws.Cells.Item(34, 11).Value = 2686.5715
$ws.Cells.Item(34, 13).Value = -2484.5715
$ws.Cells.Item(51, 8).Value = 37999
$ws.Cells.Item(51, 10).Value = 37999
$ws.Cells.Item(51, 12).Value = 37999
$ws.Cells.Item(51, 14).Value = -39471
$ws.Cells.Item(58, 8).Value = 1752
$ws.Cells.Item(58, 9).Value = 1991.3334
$ws.Cells.Item(58, 10).Value = 1321.2
$ws.Cells.Item(58, 11).Value = 1991.3334
$ws.Cells.Item(58, 12).Value = 1321.2
$ws.Cells.Item(58, 13).Value = -1788.3334
$ws.Cells.Item(58, 14).Value = -1727.2
$ws.Cells.Item(61, 8).Value = 37999
$ws.Cells.Item(61, 10).Value = 37999
$ws.Cells.Item(61, 12).Value = 37999
$ws.Cells.Item(61, 14).Value = -38695
$ws.Cells.Item(109, 8).Value = 28775.8
$ws.Cells.Item(109, 10).Value = 28775.8
$ws.Cells.Item(109, 12).Value = 28775.8
$ws.Cells.Item(109, 14).Value = -30855.8
$ws.Cells.Item(136, 8).Value = 1752
$ws.Cells.Item(136, 9).Value = 1991.3334
$ws.Cells.Item(136, 10).Value = 1321.2
$ws.Cells.Item(136, 11).Value = 5974.0002
$ws.Cells.Item(136, 12).Value = 3963.6
$ws.Cells.Item(136, 13).Value = -3424.0002
$ws.Cells.Item(136, 14).Value = -9063.6
$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Cells.Item(86, 8).Value = 803.1539
$ws.Cells.Item(86, 9).Value = 563.875
$ws.Cells.Item(86, 10).Value = 1186
$ws.Cells.Item(86, 11).Value = 1691.625
$ws.Cells.Item(86, 12).Value = 3558
$ws.Cells.Item(86, 13).Value = -505.625
$ws.Cells.Item(86, 14).Value = -5930
$ws.Cells.Item(89, 8).Value = 803.1539
$ws.Cells.Item(89, 9).Value = 563.875
$ws.Cells.Item(89, 10).Value = 1186
$ws.Cells.Item(89, 11).Value = 5074.875
$ws.Cells.Item(89, 12).Value = 10674
$ws.Cells.Item(89, 13).Value = 853.125
$ws.Cells.Item(89, 14).Value = -22530
$ws.Cells.Item(104, 8).Value = 7083.3335
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 7083.3335
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 21250.0005
$ws.Cells.Item(104, 13).Value = ""
$ws.Cells.Item(104, 14).Value = -26492.0005
$ws.Cells.Item(122, 8).Value = 3702.0212
$ws.Cells.Item(122, 9).Value = 466.66666
$ws.Cells.Item(122, 10).Value = 3922.6135
$ws.Cells.Item(122, 11).Value = 4199.99994
$ws.Cells.Item(122, 12).Value = 35303.5215
$ws.Cells.Item(122, 13).Value = -1749.99994
$ws.Cells.Item(122, 14).Value = -40203.5215
$ws.Cells.Item(132, 8).Value = 1143.4286
$ws.Cells.Item(132, 9).Value = 501
$ws.Cells.Item(132, 11).Value = 4509
$ws.Cells.Item(132, 13).Value = -1979
$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Cells.Item(57, 8).Value = 31240.666
$ws.Cells.Item(57, 10).Value = 31240.666
$ws.Cells.Item(57, 12).Value = 31240.666
$ws.Cells.Item(57, 14).Value = -32880.666
$ws.Cells.Item(80, 8).Value = 3253.4
$ws.Cells.Item(80, 9).Value = 3560.7144
$ws.Cells.Item(80, 10).Value = 2984.5
$ws.Cells.Item(80, 11).Value = 3560.7144
$ws.Cells.Item(80, 12).Value = 2984.5
$ws.Cells.Item(80, 13).Value = -2562.7144
$ws.Cells.Item(80, 14).Value = -4980.5
$ws.Cells.Item(83, 8).Value = 3253.4
$ws.Cells.Item(83, 9).Value = 3560.7144
$ws.Cells.Item(83, 10).Value = 2984.5
$ws.Cells.Item(83, 11).Value = 17803.572
$ws.Cells.Item(83, 12).Value = 14922.5
$ws.Cells.Item(83, 13).Value = -12811.572
$ws.Cells.Item(83, 14).Value = -24906.5
$ws.Cells.Item(102, 8).Value = 1446
$ws.Cells.Item(102, 9).Value = 1193
$ws.Cells.Item(102, 10).Value = 6000
$ws.Cells.Item(102, 11).Value = 1193
$ws.Cells.Item(102, 12).Value = 6000
$ws.Cells.Item(102, 13).Value = 429
$ws.Cells.Item(102, 14).Value = -9244
$ws.Cells.Item(123, 8).Value = 10210.667
$ws.Cells.Item(123, 10).Value = 10210.667
$ws.Cells.Item(123, 12).Value = 10210.667
$ws.Cells.Item(123, 14).Value = -15110.667
$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Cells.Item(22, 8).Value = 722.5599999999999
$ws.Cells.Item(22, 9).Value = 640.8333
$ws.Cells.Item(22, 10).Value = 798
$ws.Cells.Item(22, 11).Value = 640.8333
$ws.Cells.Item(22, 12).Value = 798
$ws.Cells.Item(22, 13).Value = -345.8333
$ws.Cells.Item(22, 14).Value = -1388
$ws.Cells.Item(27, 8).Value = 722.5599999999999
$ws.Cells.Item(27, 9).Value = 640.8333
$ws.Cells.Item(27, 10).Value = 798
$ws.Cells.Item(27, 11).Value = 640.8333
$ws.Cells.Item(27, 12).Value = 798
$ws.Cells.Item(27, 13).Value = -533.8333
$ws.Cells.Item(27, 14).Value = -1012
$ws.Cells.Item(141, 8).Value = 47300.5
$ws.Cells.Item(141, 9).Value = 30000
$ws.Cells.Item(141, 10).Value = 49222.777
$ws.Cells.Item(141, 11).Value = 30000
$ws.Cells.Item(141, 12).Value = 49222.777
$ws.Cells.Item(141, 14).Value = -59582.777
$ws.Cells.Item(141, 13).Value = -24820
$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Cells.Item(20, 8).Value = 1111
$ws.Cells.Item(20, 9).Value = 1111
$ws.Cells.Item(20, 11).Value = 1111
$ws.Cells.Item(20, 13).Value = -871
$ws.Cells.Item(26, 8).Value = 12507800
$ws.Cells.Item(26, 9).Value = 50000000
$ws.Cells.Item(26, 10).Value = 10400
$ws.Cells.Item(26, 11).Value = 50000000
$ws.Cells.Item(26, 12).Value = 10400
$ws.Cells.Item(26, 13).Value = -49999707
$ws.Cells.Item(26, 14).Value = -10986
$ws.Cells.Item(32, 8).Value = 19866
$ws.Cells.Item(32, 10).Value = 19866
$ws.Cells.Item(32, 12).Value = 19866
$ws.Cells.Item(32, 14).Value = -20500
$ws.Cells.Item(39, 8).Value = 17487.25
$ws.Cells.Item(39, 10).Value = 17487.25
$ws.Cells.Item(39, 12).Value = 17487.25
$ws.Cells.Item(39, 14).Value = -18313.25
$ws.Cells.Item(51, 8).Value = 4100
$ws.Cells.Item(51, 9).Value = 4100
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 4100
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -3590
$ws.Cells.Item(51, 14).Value = ""
$ws.Cells.Item(96, 8).Value = 2069.182
$ws.Cells.Item(96, 9).Value = 1651
$ws.Cells.Item(96, 10).Value = 2226
$ws.Cells.Item(96, 11).Value = 1651
$ws.Cells.Item(96, 12).Value = 2226
$ws.Cells.Item(96, 13).Value = -278
$ws.Cells.Item(96, 14).Value = -4972
